$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right below
#    the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Big Santa Fortune for Free - Review
#    of Features, RTP, and Payout") right before the final (DALLE prompt)
#    paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$origStart = $lastPara.Range.Start
$insertionPoint = $d.Range($origStart, $origStart)

$boldText = "Play Big Santa Fortune for Free – Review of Features, RTP, and Payout"

$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($frag)

# Split the paragraph right after the bold text we just typed in, so the
# DALLE-prompt paragraph that originally started at $origStart becomes its
# own paragraph again.
$splitPoint = $origStart + $boldText.Length
$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3) Replace the DALLE image-prompt text with the new meta-description text
#    (keeping the paragraph's italic formatting untouched).
# ---------------------------------------------------------------------------
$newLastIdx = $d.Paragraphs.Count
$dallePara = $d.Paragraphs($newLastIdx)

$oldText = 'DALLE, please create a cartoon-style feature image for the game "Big Santa Fortune" that features a happy Maya warrior with glasses. The image should be eye-catching and engaging, with bright colors and the Maya warrior holding a fishing rod, standing next to Santa Claus, who is holding the Big Sticky Fish. The background should include winter and summer elements, such as snow and sand, and include Christmas-themed symbols like presents and Christmas balls. Overall, the image should convey the fun and festive nature of the game while highlighting the important elements, such as the bonus feature and the potential for big wins.'
$newText = 'Read our review of Big Santa Fortune slot, a unique and entertaining online game that combines Christmas and fishing themes. Play it for free and discover its features, RTP, and payout!'

$find = $dallePara.Range.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
